# Adds two new donor sign-ups to "Donors" and one new patient request to
# "Patients", and clears a stray date-number-format that had been applied
# to Patients!J16 (its value is a plain number again, matching the rest
# of column J).

$wb = $excel.ActiveWorkbook
$wsDonors   = $wb.Worksheets.Item("Donors")
$wsPatients = $wb.Worksheets.Item("Patients")

# ---------------------------------------------------------------------
# Donors (sheet1): two new rows, 20 and 21
# ---------------------------------------------------------------------
$wsDonors.Range("A20").Value = "69283187b147c0f70a38e0ea"
$wsDonors.Range("B20").Value = "Ashwini Shenoy B"
$wsDonors.Range("C20").Value = "ashenoyb@gmail.com"
# Phone numbers are all-digit strings that must stay text (matches the
# sheet's numberStoredAsText ignored-error range) rather than becoming a
# numeric value - force text entry, then drop the resulting "Text" format
# so no explicit style lingers on the cell (same bare <c t="str"> shape as
# every other phone-number cell in the column).
$wsDonors.Range("D20").NumberFormat = "@"
$wsDonors.Range("D20").Value = "7026438371"
$wsDonors.Range("D20").ClearFormats()
$wsDonors.Range("E20").Value = "A+"
$wsDonors.Range("F20").Value = "My Current Location"
$wsDonors.Range("G20").Value = 77.5483533
$wsDonors.Range("H20").Value = 12.9565222
$wsDonors.Range("I20").Value = $false
$wsDonors.Range("J20").Value = $false
$wsDonors.Range("K20").Value = $true
$wsDonors.Range("L20").Value = 45988.694550451386

$wsDonors.Range("A21").Value = "692834afb147c0f70a38e12a"
$wsDonors.Range("B21").Value = "Ashwini Shenoy B"
$wsDonors.Range("C21").Value = "shenoybashwini@gmail.com"
$wsDonors.Range("D21").NumberFormat = "@"
$wsDonors.Range("D21").Value = "1234567890"
$wsDonors.Range("D21").ClearFormats()
$wsDonors.Range("E21").Value = "O+"
$wsDonors.Range("F21").Value = "My Current Location"
$wsDonors.Range("G21").Value = 77.5483533
$wsDonors.Range("H21").Value = 12.9565222
$wsDonors.Range("I21").Value = $false
$wsDonors.Range("J21").Value = $false
$wsDonors.Range("K21").Value = $true
$wsDonors.Range("L21").Value = 45988.703903194444

# L21 carries the same date number-format already used by an existing
# timestamp cell (Patients!J16) - copy formats only (before J16's own
# format gets cleared below), so the value just written above is left
# untouched.
$wsPatients.Range("J16").Copy() | Out-Null
$wsDonors.Range("L21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Patients (sheet2): clear the odd date format on J16, then add row 17
# ---------------------------------------------------------------------
$wsPatients.Range("J16").ClearFormats()

$wsPatients.Range("A17").Value = "692831b6b147c0f70a38e0f1"
$wsPatients.Range("B17").Value = "Ashwini Shenoy B"
$wsPatients.Range("C17").Value = "ashwinishenoyb@gmail.com"
$wsPatients.Range("D17").NumberFormat = "@"
$wsPatients.Range("D17").Value = "7026438371"
$wsPatients.Range("D17").ClearFormats()
$wsPatients.Range("E17").Value = "A+"
$wsPatients.Range("F17").Value = "My Current Location"
$wsPatients.Range("G17").Value = 77.5483533
$wsPatients.Range("H17").Value = 12.9565222
$wsPatients.Range("I17").Value = "High"
$wsPatients.Range("J17").Value = 45988.695100902776
